$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.3048080303191223, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 5.014808316549482)
    3  = @(0.6753301551942219, 114.8270160096505, 337.1190423067083, 8.660232485948974, 461.281620957502)
    4  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    5  = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 8.660232485948974, 38.00504557982321)
    6  = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 31.61296591696135)
    7  = @(0.0008583669626518464, 0.04240448674262143, 3.900430680208489, 0.496779210170732, 4.440472744084493)
    8  = @(0.127881588408715, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 11.26139829420982)
    9  = @(0.127881588408715, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.449980674824537)
    10 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    11 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
